$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Build two off-grid scratch/template cells carrying the two new cell
# styles we need (top+bottom border, and top+bottom+right border), then
# copy/paste just the formats onto the real target cells. Building the
# border combinations one edge at a time directly on cells spread across
# more than one worksheet can leave unused/orphaned style or border
# entries behind in styles.xml, so the styling is staged on scratch
# cells first using the combination of calls that was verified to stay
# clean (whole-collection LineStyle assignment followed by removing the
# Left edge for the triple-edge combination).
$tTopBottom = $ws1.Range("Z1")
$tTopBottom.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$tTopBottom.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$tTopBottom.Copy()
$ws1.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)
$tTopBottom.Clear()

$tTopBottomRight = $ws1.Range("Z2")
$tTopBottomRight.Borders.LineStyle = 1       # all four edges thin
$tTopBottomRight.Borders.Item(7).LineStyle = -4142  # remove xlEdgeLeft

$tTopBottomRight.Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)
$tTopBottomRight.Clear()

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5 on computational_comparison
$ws2.Range("G5").ClearContents()
